# Refresh the cryptocurrency price table (Price/Volume columns, plus the
# coin rows shifted down one slot by the newly-inserted "WrappedeETH" row).
# Generated from the upstream GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (row, column, new text value). Column 4 is "Price" (D) and
# holds values that look numeric (e.g. "595.80", "67.697.92"); Excel's COM
# Value setter auto-converts those to numbers/doubles, which would lose the
# source formatting (and the two-dot "thousands" values aren't even valid
# numbers). Force those cells to text first, then clear the temporary "@"
# number format back to the sheet's default (General, no explicit style)
# so the cell's style index is left untouched, matching the original file.
$updates = @(
    @(2, 4, "67.697.92"),
    @(2, 5, "  -0.65%  "),
    @(3, 4, "3.790.79"),
    @(3, 5, "  +0.54%  "),
    @(4, 5, "  +0.00%  "),
    @(5, 4, "595.80"),
    @(5, 5, "  +0.42%  "),
    @(6, 4, "166.90"),
    @(6, 5, "  -0.25%  "),
    @(7, 4, "3.790.46"),
    @(7, 5, "  +0.59%  "),
    @(8, 5, "  -0.01%  "),
    @(9, 5, "  +0.24%  "),
    @(10, 5, "  +0.11%  "),
    @(11, 5, "  -0.97%  "),
    @(12, 5, "  -0.06%  "),
    @(13, 5, "  -1.93%  "),
    @(14, 4, "36.06"),
    @(14, 5, "  +0.04%  "),
    @(15, 4, "4.426.28"),
    @(15, 5, "  +0.56%  "),
    @(16, 4, "3.779.98"),
    @(16, 5, "  +0.78%  "),
    @(17, 4, "18.54"),
    @(17, 5, "  +4.02%  "),
    @(18, 4, "67.644.70"),
    @(19, 5, "  +1.01%  "),
    @(20, 5, "  +0.12%  "),
    @(21, 4, "10.00"),
    @(21, 5, "  -7.23%  "),
    @(22, 4, "459.33"),
    @(22, 5, "  -0.91%  "),
    @(23, 4, "0.698"),
    @(23, 5, "  +0.31%  "),
    @(24, 5, "  +4.27%  "),
    @(25, 5, "  -0.37%  "),
    @(26, 4, "12.16"),
    @(26, 5, "  +2.95%  "),
    @(27, 5, "  -3.05%  "),
    @(28, 5, "  -0.01%  "),
    @(29, 4, "9.99"),
    @(29, 5, "  -0.90%  "),
    @(30, 2, "WrappedeETH"),
    @(30, 3, "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"),
    @(30, 4, "3.940.51"),
    @(30, 5, "  +0.54%  "),
    @(31, 2, "PancakeSwap"),
    @(31, 3, "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"),
    @(31, 4, "2.77"),
    @(31, 5, "  -0.17%  "),
    @(32, 2, "ImmutableX"),
    @(32, 3, "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @(32, 4, "2.25"),
    @(32, 5, "  +4.50%  "),
    @(33, 2, "NEARProtocol"),
    @(33, 3, "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"),
    @(33, 4, "7.21"),
    @(33, 5, "  -1.01%  "),
    @(34, 2, "EthereumClassic"),
    @(34, 3, "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"),
    @(34, 4, "29.65"),
    @(34, 5, "  -0.60%  "),
    @(35, 2, "Binance-PegBSC-USD"),
    @(35, 3, "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"),
    @(35, 4, "1.00"),
    @(35, 5, "  +0.00%  "),
    @(36, 2, "Aptos"),
    @(36, 3, "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"),
    @(36, 4, "9.06"),
    @(36, 5, "  -0.52%  "),
    @(37, 2, "Hedera"),
    @(37, 3, "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @(37, 4, "0.0999"),
    @(37, 5, "  -0.42%  "),
    @(38, 2, "dogwifhat"),
    @(38, 3, "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"),
    @(38, 4, "3.36"),
    @(38, 5, "  -1.90%  "),
    @(39, 2, "Kaspa"),
    @(39, 3, "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @(39, 4, "0.138"),
    @(39, 5, "  -0.30%  "),
    @(40, 2, "Mantle"),
    @(40, 3, "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"),
    @(40, 4, "0.994"),
    @(40, 5, "  -0.65%  "),
    @(41, 2, "Filecoin"),
    @(41, 3, "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"),
    @(41, 4, "5.77"),
    @(41, 5, "  +0.25%  "),
    @(42, 2, "FirstDigitalUSD"),
    @(42, 3, "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"),
    @(42, 4, "0.999"),
    @(43, 2, "USDe"),
    @(43, 3, "https://coinranking.com/coin/exbfr2U-0+usde-usde"),
    @(43, 4, "1.00"),
    @(43, 5, "  +0.01%  "),
    @(44, 2, "OKB"),
    @(44, 3, "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"),
    @(44, 4, "48.03"),
    @(44, 5, "  +2.42%  "),
    @(45, 2, "Arweave"),
    @(45, 3, "https://coinranking.com/coin/7XWg41D1+arweave-ar"),
    @(45, 4, "43.80"),
    @(45, 5, "  -1.18%  "),
    @(46, 2, "TheGraph"),
    @(46, 3, "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"),
    @(46, 4, "0.298"),
    @(46, 5, "  -0.56%  "),
    @(47, 2, "Monero"),
    @(47, 3, "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @(47, 4, "149.87"),
    @(47, 5, "  +3.08%  "),
    @(48, 2, "Cosmos"),
    @(48, 3, "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"),
    @(48, 4, "8.27"),
    @(48, 5, "  -1.42%  "),
    @(49, 2, "EnergySwap"),
    @(49, 3, "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @(49, 4, "27.04"),
    @(49, 5, "  +7.65%  "),
    @(50, 2, "Bittensor"),
    @(50, 3, "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"),
    @(50, 4, "390.01"),
    @(50, 5, "  +0.55%  "),
    @(51, 2, "Stacks"),
    @(51, 3, "https://coinranking.com/coin/mMPrMcB7+stacks-stx"),
    @(51, 4, "1.82"),
    @(51, 5, "  -4.26%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $value = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    if ($col -eq 4) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}
